# "Stor ändring av lager med koppling till firebase"
# Update the "Nya enheter" (New units) sheet: the first product line item
# ("Iphone 16" / "Tre olika färger") is replaced by four rows of
# "Iphone 15 Pro" in different colour variants, each with quantity 5, and
# the following two existing products (Samsung 24 Ultra, Macbook Air 13)
# are pushed down two rows, keeping their original quantities.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nya enheter")

# Row 3: Iphone 15 Pro - Svart XX
$ws.Range("C3").Value = "Iphone 15 Pro"
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = "Svart XX"

# Row 4: Iphone 15 Pro - Svart XY
$ws.Range("C4").Value = "Iphone 15 Pro"
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = "Svart XY"

# Row 5: Iphone 15 Pro - Vitt YY
$ws.Range("C5").Value = "Iphone 15 Pro"
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = "Vitt YY"

# Row 6: Iphone 15 Pro - Guld YX (new row, previously only had a code)
$ws.Range("C6").Value = "Iphone 15 Pro"
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = "Guld YX"

# Row 7: Samsung 24 Ultra (moved down from row 4)
$ws.Range("C7").Value = "Samsung 24 Ultra"
$ws.Range("D7").Value = 3

# Row 8: Macbook Air 13 (moved down from row 5)
$ws.Range("C8").Value = "Macbook Air 13"
$ws.Range("D8").Value = 2

# Update the active selection to reflect where the user ended up editing
$ws.Range("E11").Select()
